$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 = R5: value fix from 4K7 to 100 ohm, and JLCPCB part # fix
$ws.Range("A22").Value = 100
$ws.Range("D22").Value = "C22775"

# The "22775" portion of the new part number keeps the pasted-in 宋体 run formatting
$ws.Range("D22").Characters(2, 5).Font.Name = "宋体"
$ws.Range("D22").Characters(2, 5).Font.Size = 11
$ws.Range("D22").Characters(2, 5).Font.Color = 0

$ws.Range("D23").Select()
